# Update countries & provincias Spain
# - Refresh the "last updated" timestamp
# - Refresh case numbers for several countries
# - Bangladesh's case count overtakes Eslovaquia's, so it moves up one
#   rank in the (descending, by total cases) table; every country that
#   was between Bangladesh's old spot and its new spot shifts down one row

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Last updated timestamp (row 1) ---
$ws.Range("A1").Value = "Datos actualizados a 13 de Abril de 2020 a las 10:52"

# --- Austria (row 20) ---
$ws.Range("D20").Value = 7343
$ws.Range("E20").Value = 6251
$ws.Range("F20").Value = 239
$ws.Range("G20").Value = 18
$ws.Range("H20").Value = 368

# --- Eslovenia (row 68) ---
$ws.Range("B68").Value = 1212
$ws.Range("C68").Value = 7
$ws.Range("D68").Value = 152
$ws.Range("E68").Value = 1005
$ws.Range("F68").Value = 34
$ws.Range("G68").Value = 2
$ws.Range("H68").Value = 55

# --- Hong Kong (row 74) ---
$ws.Range("B74").Value = 1010
$ws.Range("C74").Value = 5
$ws.Range("D74").Value = 397
$ws.Range("E74").Value = 609

# --- Kazajistan (row 75) ---
$ws.Range("B75").Value = 979
$ws.Range("C75").Value = 28
$ws.Range("D75").Value = 110
$ws.Range("E75").Value = 857
$ws.Range("G75").Value = 2
$ws.Range("H75").Value = 12

# --- Rows 79-89: Bangladesh jumps ahead of Eslovaquia, pushing the
#     countries formerly ranked 79-88 down to 80-89 ---

# row 79 -> now Banglades, with refreshed data
$ws.Range("A79").Value = "Banglades"
$ws.Range("B79").Value = 803
$ws.Range("C79").Value = 182
$ws.Range("D79").Value = 42
$ws.Range("E79").Value = 722
$ws.Range("F79").Value = 1
$ws.Range("G79").Value = 5
$ws.Range("H79").Value = 39

# row 80 -> now Eslovaquia (former row 79 data)
$ws.Range("A80").Value = "Eslovaquia"
$ws.Range("B80").Value = 742
$ws.Range("C80").Value = 0
$ws.Range("D80").Value = 23
$ws.Range("E80").Value = 717
$ws.Range("F80").Value = 5
$ws.Range("G80").Value = 0
$ws.Range("H80").Value = 2

# row 81 -> now Oman (former row 80 data)
$ws.Range("A81").Value = "Oman"
$ws.Range("B81").Value = 727
$ws.Range("C81").Value = 128
$ws.Range("D81").Value = 124
$ws.Range("E81").Value = 599
$ws.Range("F81").Value = 3
$ws.Range("G81").Value = 0
$ws.Range("H81").Value = 4

# row 82 -> now Crucero (former row 81 data)
$ws.Range("A82").Value = "Crucero"
$ws.Range("B82").Value = 712
$ws.Range("C82").Value = 0
$ws.Range("D82").Value = 619
$ws.Range("E82").Value = 82
$ws.Range("F82").Value = 10
$ws.Range("G82").Value = 0
$ws.Range("H82").Value = 11

# row 83 -> now Tunez (former row 82 data)
$ws.Range("A83").Value = "Tunez"
$ws.Range("B83").Value = 707
$ws.Range("C83").Value = 0
$ws.Range("D83").Value = 43
$ws.Range("E83").Value = 633
$ws.Range("F83").Value = 85
$ws.Range("G83").Value = 0
$ws.Range("H83").Value = 31

# row 84 -> now Bulgaria (former row 83 data)
$ws.Range("A84").Value = "Bulgaria"
$ws.Range("B84").Value = 676
$ws.Range("C84").Value = 1
$ws.Range("D84").Value = 71
$ws.Range("E84").Value = 574
$ws.Range("F84").Value = 36
$ws.Range("G84").Value = 2
$ws.Range("H84").Value = 31

# row 85 -> now Cuba (former row 84 data)
$ws.Range("A85").Value = "Cuba"
$ws.Range("B85").Value = 669
$ws.Range("C85").Value = 0
$ws.Range("D85").Value = 92
$ws.Range("E85").Value = 559
$ws.Range("F85").Value = 11
$ws.Range("G85").Value = 0
$ws.Range("H85").Value = 18

# row 86 -> now Letonia (former row 85 data)
$ws.Range("A86").Value = "Letonia"
$ws.Range("B86").Value = 653
$ws.Range("C86").Value = 2
$ws.Range("D86").Value = 16
$ws.Range("E86").Value = 632
$ws.Range("F86").Value = 2
$ws.Range("G86").Value = 0
$ws.Range("H86").Value = 5

# row 87 -> now Principado de Andorra (former row 86 data)
$ws.Range("A87").Value = "Principado de Andorra"
$ws.Range("B87").Value = 638
$ws.Range("C87").Value = 0
$ws.Range("D87").Value = 128
$ws.Range("E87").Value = 481
$ws.Range("F87").Value = 17
$ws.Range("G87").Value = 0
$ws.Range("H87").Value = 29

# row 88 -> now Republica de Chipre (former row 87 data)
$ws.Range("A88").Value = "Republica de Chipre"
$ws.Range("B88").Value = 633
$ws.Range("C88").Value = 0
$ws.Range("D88").Value = 65
$ws.Range("E88").Value = 557
$ws.Range("F88").Value = 8
$ws.Range("G88").Value = 0
$ws.Range("H88").Value = 11

# row 89 -> now Libano (former row 88 data)
$ws.Range("A89").Value = "Libano"
$ws.Range("B89").Value = 630
$ws.Range("C89").Value = 0
$ws.Range("D89").Value = 80
$ws.Range("E89").Value = 530
$ws.Range("F89").Value = 34
$ws.Range("G89").Value = 0
$ws.Range("H89").Value = 20

# row 90 (Afganistan) is unchanged - Bangladesh's new total (803) still
# trails Libano's former total (630) only down to Afganistan's 607, so
# Afganistan keeps its row.

# --- Vietnam (row 112) ---
$ws.Range("D112").Value = 145
$ws.Range("E112").Value = 117
